$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D13").Value = "https://stats.ifp.uni-mainz.de/ba-ccs-track/befragung-rewb.html"
$ws.Range("F13").Value = "exercises/e12.html"
$ws.Range("F14").Select()
